# Applies the text edits from the commit to the single slide of the deck:
#  - "Diamond 6"  shape: "Reached last root node in DTree?"
#       -> "Exhausted all root node in DTree?"
#       (the leading run is split into two runs with identical formatting:
#        "Exhausted all root " + "node in ")
#  - "Diamond 10" shape: "Reached last category for the selected root node?"
#       -> "Exhausted all category for the selected root node?"
#       (split into "Exhausted all category " + "for the selected root node?")
#
# Each replacement is done as two separate Characters(...).Text assignments
# (rather than one assignment over the whole original run) so PowerPoint
# keeps them as two distinct runs with identical, untouched <a:rPr/>
# formatting - matching how the real authoring tool produced the edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Shape "Diamond 6" : "Reached last root node in " + "DTree" + "?" ----
$sh1 = $s.Shapes.Item("Diamond 6")
$tr1 = $sh1.TextFrame.TextRange
Write-Host ("Diamond 6 before: [" + $tr1.Text + "]")

# The original leading run "Reached last root node in " is 26 characters.
# Replace its first 19 characters with the new first half of the text...
$newFirstHalf = "Exhausted all root "
$partA = $tr1.Characters(1, 19)
$partA.Text = $newFirstHalf

# ...then replace the remaining 7 characters of that original run (now
# shifted by the new first half's length) with the new second half.
$newSecondHalf = "node in "
$tr1Again = $sh1.TextFrame.TextRange
$partB = $tr1Again.Characters($newFirstHalf.Length + 1, 26 - 19)
$partB.Text = $newSecondHalf

Write-Host ("Diamond 6 after: [" + $sh1.TextFrame.TextRange.Text + "]")

# ---- Shape "Diamond 10" : "Reached last category for the selected root node?" ----
$sh2 = $s.Shapes.Item("Diamond 10")
$tr2 = $sh2.TextFrame.TextRange
Write-Host ("Diamond 10 before: [" + $tr2.Text + "]")

# The whole text is a single 49-character run. Replace the first 22
# characters with the new first half...
$newFirstHalf2 = "Exhausted all category "
$part2A = $tr2.Characters(1, 22)
$part2A.Text = $newFirstHalf2

# ...then replace the rest of the original run with the new second half.
$newSecondHalf2 = "for the selected root node?"
$tr2Again = $sh2.TextFrame.TextRange
$part2B = $tr2Again.Characters($newFirstHalf2.Length + 1, $tr2Again.Length - $newFirstHalf2.Length)
$part2B.Text = $newSecondHalf2

Write-Host ("Diamond 10 after: [" + $sh2.TextFrame.TextRange.Text + "]")
